# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded for the
# 2d599492-60f7-42e2-80b5-c099bbae70cf.md file (row 3 on each sheet),
# as would happen after a fresh handback report was generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2017-02-21 09:07:54"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (L) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2017-02-21 09:07:36"
$wsZhCn.Range("L3").Value = "2017-02-21 09:08:35"

# --- de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (L) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2017-02-21 09:07:54"
$wsDeDe.Range("L3").Value = "2017-02-21 09:08:58"
